$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "data as of" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Junio de 2020 a las 14:26"

# Benin moved up in the ranking (now ranked just behind Hong Kong instead of
# behind Cabo Verde); update the three affected country-name cells so the rows
# keep matching their refreshed case counts below.
$ws.Range("A125").Value = "Benin"
$ws.Range("A126").Value = "Tunez"
$ws.Range("A127").Value = "Cabo Verde"

# Refreshed COVID-19 case numbers per country (columns B:H = Casos totales,
# Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
# Row 4
$ws.Range("B4").Value = 2637241
$ws.Range("C4").Value = 164
$ws.Range("D4").Value = 1093545
$ws.Range("E4").Value = 1415258

# Row 5
$ws.Range("B5").Value = 1345470
$ws.Range("C5").Value = 216
$ws.Range("E5").Value = 553963
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 57659

# Row 7
$ws.Range("B7").Value = 549991
$ws.Range("C7").Value = 794
$ws.Range("D7").Value = 322364
$ws.Range("E7").Value = 211123

# Row 28
$ws.Range("B28").Value = 61790
$ws.Range("C28").Value = 315
$ws.Range("D28").Value = 45213
$ws.Range("E28").Value = 16190
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 387

# Row 33
$ws.Range("B33").Value = 50223
$ws.Range("C33").Value = 76
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 6107

# Row 35
$ws.Range("B35").Value = 47151
$ws.Range("C35").Value = 1749
$ws.Range("D35").Value = 22974
$ws.Range("E35").Value = 22338
$ws.Range("G35").Value = 83
$ws.Range("H35").Value = 1839

# Row 36
$ws.Range("B36").Value = 45524
$ws.Range("C36").Value = 582
$ws.Range("D36").Value = 36313
$ws.Range("E36").Value = 8861
$ws.Range("G36").Value = 2
$ws.Range("H36").Value = 350

# Row 61
$ws.Range("D61").Value = 9229
$ws.Range("E61").Value = 6488
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 533

# Row 66
$ws.Range("B66").Value = 12751
$ws.Range("C66").Value = 76
$ws.Range("D66").Value = 11612
$ws.Range("E66").Value = 534
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = 605

# Row 74
$ws.Range("B74").Value = 8199
$ws.Range("C74").Value = 251
$ws.Range("D74").Value = 5426
$ws.Range("E74").Value = 2751

# Row 77
$ws.Range("B77").Value = 6939
$ws.Range("C77").Value = 112
$ws.Range("D77").Value = 1050
$ws.Range("E77").Value = 5722
$ws.Range("G77").Value = 10
$ws.Range("H77").Value = 167

# Row 101
$ws.Range("B101").Value = 2725
$ws.Range("C101").Value = 34
$ws.Range("D101").Value = 2155
$ws.Range("E101").Value = 463

# Row 117
$ws.Range("B117").Value = 1745
$ws.Range("C117").Value = 5
$ws.Range("D117").Value = 1170
$ws.Range("E117").Value = 541

# Row 125
$ws.Range("B125").Value = 1187
$ws.Range("C125").Value = 38
$ws.Range("D125").Value = 323
$ws.Range("E125").Value = 845
$ws.Range("G125").Value = 3
$ws.Range("H125").Value = 19

# Row 126
$ws.Range("B126").Value = 1169
$ws.Range("D126").Value = 1029
$ws.Range("E126").Value = 90
$ws.Range("H126").Value = 50

# Row 127
$ws.Range("B127").Value = 1155
$ws.Range("D127").Value = 570
$ws.Range("E127").Value = 573
$ws.Range("H127").Value = 12

# Row 163
$ws.Range("D163").Value = 93
$ws.Range("E163").Value = 163

Write-Host "Updated countries & provincias Spain data"